$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)
$ws.Range("A1").NumberFormat = "@"
$ws.Range("A1").Value = "3103040004"
$ws.Range("A1").NumberFormat = "General"
Write-Host "done"
